# FinancieroG.xlsx - "Mejoras en presentacion de mensajes y paginacion"
#
# 1) Sheet "Proveedores": normalize the date-format style used in column B
#    (rows 2-245) from the old numFmt (style index 2 / s=3 on the last row)
#    to the new numFmt (style index 4), then append 8 new transaction rows
#    (246-253), the last one using a date-only style (style index 5).
# 2) Sheet "Resumen": append 2 new summary rows (11-12) for the new
#    "Prueba" and "Pepito" providers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Proveedores"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Proveedores")

# Re-stamp the existing date cells (B2:B245) with the new (equivalent)
# datetime number format so they all share style index 4.
$ws.Range("B2:B245").NumberFormat = "yyyy-mm-dd h:mm:ss"

# New rows appended at the bottom of the table.
$newRows = @(
    @{ Row = 246; A = 245; B = 45924; C = "Prueba"; D = "Factura"; E = "Factura"; F = 10000; G = "F-001"; H = "Inactiva" },
    @{ Row = 247; A = 246; B = 45924; C = "Prueba"; D = "Abono";   E = "Abono";   F = 3000;  G = "F-001"; H = "Inactiva" },
    @{ Row = 248; A = 247; B = 45924; C = "Prueba"; D = "Factura"; E = "Factura - la factura se realizo por : 3000  + saldo anterior 7000 = total 10000"; F = 10000; G = "F-002"; H = "Inactiva" },
    @{ Row = 249; A = 248; B = 45924; C = "Prueba"; D = "Abono";   E = "Abono";   F = 5000;  G = "F-002"; H = "Inactiva" },
    @{ Row = 250; A = 249; B = 45924; C = "Prueba"; D = "Abono";   E = "Abono";   F = 3500;  G = "F-002"; H = "Inactiva" },
    @{ Row = 251; A = 250; B = 45924; C = "Prueba"; D = "Factura"; E = "Factura - la factura se realizo por : 4000  + saldo anterior 1500 = total 5500"; F = 5500; G = "F-003"; H = "Activa" },
    @{ Row = 252; A = 251; B = 45924; C = "Prueba"; D = "Abono";   E = "Abono";   F = 2500;  G = "F-003"; H = "Activa" },
    @{ Row = 253; A = 252; B = 45924; C = "Pepito"; D = "Factura"; E = "Factura"; F = 50000; G = "F-001"; H = "Activa" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    if ($rowNum -eq 253) {
        $ws.Cells.Item($rowNum, 2).NumberFormat = "yyyy-mm-dd"
    } else {
        $ws.Cells.Item($rowNum, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"
    }
}

# ---------------------------------------------------------------------
# Sheet "Resumen"
# ---------------------------------------------------------------------
$wsResumen = $wb.Worksheets.Item("Resumen")

$wsResumen.Cells.Item(11, 1).Value = 10
$wsResumen.Cells.Item(11, 2).Value = "Prueba"
$wsResumen.Cells.Item(11, 3).Value = 5500
$wsResumen.Cells.Item(11, 4).Value = 2500
$wsResumen.Cells.Item(11, 5).Value = 3000

$wsResumen.Cells.Item(12, 1).Value = 11
$wsResumen.Cells.Item(12, 2).Value = "Pepito"
$wsResumen.Cells.Item(12, 3).Value = 50000
$wsResumen.Cells.Item(12, 4).Value = 0
$wsResumen.Cells.Item(12, 5).Value = 50000
